$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.37%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.016"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.01%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.057"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-20.80%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.803"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.21%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.777"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.30%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9180"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.87%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1742"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.99%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07834"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08736"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.41%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03097"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.18%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001523"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.39%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005930"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.55%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.460"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.64%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.267"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.82%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.32%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.188"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.49%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.08%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.13%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.01%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004472"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.35%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001250"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.49%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01738"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.79%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04757"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007343"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.77%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1356"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.37%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002159"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.11%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01072"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.92%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006064"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-7.54%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.31%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003511"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-59.74%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8235"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.35%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.31%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.31%"
